# The workbook tracks a rolling 11-month window of daily totals (columns
# B:AF = day-of-month 1..31, rows 2.. = one row per month, column A = month
# name label). This edit rolls the window forward: the oldest month (Jun)
# is dropped from the top and three new months (Apr, Jun, Jul) are appended
# at the bottom, while every other month's data is refreshed with newer
# figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Month labels for rows 2..12 (column A), in order.
$months = @("Aug","Sep","Oct","Nov","Dec","Jan","Feb","Mar","Apr","Jun","Jul")

# Daily totals for rows 2..12, columns B..AF (31 values per row).
$data = @(
    @(0,0,0,0,0,0,0,1148.89,1253.01,1430.72,2568.03,2782.93,1676.38,1334.57,1027.22,955.45,969.25,968.39,1053.27,592.63,677.52,905.17,855.46,753.42,1478.55,1740.64,1114.77,742,489.75,526.1,991.3),
    @(1483.91,1879.41,1353.95,493.89,474.77,443.32,548.55,1054.67,1199.69,1025.38,891.25,1038.95,1173.73,1230.03,2251.04,1905.21,764.72,923.21,1046.92,786.12,931.89,1028.86,798.26,564.65,465.85,469.15,525,709.85,1244.5,1452.35,0),
    @(348,249,249,82,803,2397,2316,1269,119,119,119,1137,2509,2509,1586,374,478,642,1732,2507,2621,1606,762,123,591,989,1305,1305,704,82,164),
    @(123,186,324,406,156,82,82,76,99,123,123,111,70,35,0,75.55,91.2,283,207.2,149.2,200,71,149.57,128.67,172.87,113.87,76.16,37,71,74.7,0),
    @(136.42,132,39,37,68,68,71.4,132,132,78,160.7,160.7,123.7,39,0,0,0,64.6,64.6,39.26,78.25,86.25,47,39,74,74,37,0,47,47,75.6),
    @(73.6,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,199,199,159,0,111,0,0,0,0,0,0,0,0),
    @(0,0,0,0,0,0,0,34,37.4,37.4,0,0,0,0,149,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0),
    @(0,0,149,129,189,243,207,262,262,187,191,550,550,715,930,163,49,0,0,0,0,0,0,0,0,0,0,0,0,52,41),
    @(198,198,237,104.45,1303.85,1567.96,1160.96,986.2,799.2,308,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0),
    @(0,0,0,0,0,0,0,0,0,0,0,43,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0),
    @(0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,179,219,219,179,0,0,0)
)

$rowCount = $data.Count
$colCount = $data[0].Count

# Build a 2-D array for a single bulk write into B2:AF12.
$arr = New-Object 'object[,]' $rowCount,$colCount
for ($i = 0; $i -lt $rowCount; $i++) {
    for ($j = 0; $j -lt $colCount; $j++) {
        $arr[$i,$j] = $data[$i][$j]
    }
}
$ws.Range("B2:AF12").Value = $arr

# Write the month labels in column A for rows 2..12.
for ($i = 0; $i -lt $months.Count; $i++) {
    $r = 2 + $i
    $ws.Cells.Item($r, 1).Value = $months[$i]
}

# Rows 10-12 are brand new; give their column-A label cells the same
# formatting (bold, bordered, centered) as the existing month labels by
# copying the format from row 9's label cell.
$ws.Range("A9").Copy() | Out-Null
$ws.Range("A10:A12").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
